# Update column N (ratio) values on sheet Tab23 per the source diff
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab23")

$ws.Range("N13").Value = 2525.2400042378199
$ws.Range("N23").Value = 2373.8133079182198
$ws.Range("N38").Value = 2258.1855425952399
$ws.Range("N45").Value = 3336.2071632669099
$ws.Range("N61").Value = 2590.5500867576102
$ws.Range("N62").Value = 2597.67579128335
$ws.Range("N63").Value = 2973.5091189032701
$ws.Range("N64").Value = 3047.2970802363898
$ws.Range("N65").Value = 2804.1083144570098
$ws.Range("N66").Value = 2916.5117556755899
$ws.Range("N67").Value = 2511.50676411314
$ws.Range("N68").Value = 2755.4966645743498
$ws.Range("N69").Value = 2230.1906916993898
$ws.Range("N70").Value = 2320.2769207746501
$ws.Range("N71").Value = 2590.5500867576102
$ws.Range("N72").Value = 2262.22687253816
$ws.Range("N73").Value = 2437.6284547318101
$ws.Range("N74").Value = 3352.0306002006701
$ws.Range("N75").Value = 2832.1198234047301
$ws.Range("N76").Value = 3051.80122409674
$ws.Range("N77").Value = 3450.4861828215398
$ws.Range("N78").Value = 2565.1844073995699
$ws.Range("N79").Value = 3002.9781012830599
$ws.Range("N80").Value = 2612.2435405891802
$ws.Range("N81").Value = 2971.0449106783599
$ws.Range("N82").Value = 2357.5412526099699
$ws.Range("N83").Value = 2048.36264164764
$ws.Range("N84").Value = 2704.5365449258402
$ws.Range("N85").Value = 2560.2928867453902
$ws.Range("N86").Value = 2938.7327190576202
$ws.Range("N87").Value = 3136.9501503757401
$ws.Range("N88").Value = 3430.63190255068
$ws.Range("N89").Value = 2365.7601031254299
$ws.Range("N90").Value = 2515.0144422708399
$ws.Range("N91").Value = 2586.6169510907398
$ws.Range("N92").Value = 2782.3243861455298
$ws.Range("N93").Value = 2343.3037766356401
$ws.Range("N94").Value = 2683.4725649816901
$ws.Range("N95").Value = 2383.6527428618601
$ws.Range("N96").Value = 2536.0402697965901
$ws.Range("N97").Value = 2235.95057650024

